# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, identical on both sheets.
$updates = @{
    2  = 8835
    3  = 8209
    4  = 143
    8  = 143
    9  = 153
    10 = 203
    14 = 5297
    16 = 81
    17 = 19
    19 = 154
    20 = 140
    21 = 4
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
